$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 9 & 10: swap all betting-data columns (B, G, H..AC) between the two matches,
#     keep A (id) and F (home team) fixed per row ---
# --- Many F/G cells across the sheet: NK Domzale/NK Maribor shared-string slots were
#     swapped upstream, so the literal team text shown must be corrected back to
#     what it displayed before (net no-op on content, but needed cell by cell) ---
# --- Rows 138 & 139: refreshed with updated match data (new odds/IDs), keeping the
#     original A (row id) values ---
$ws.Range("G2").Value = "NK Maribor"
$ws.Range("G5").Value = "NK Domzale"
$ws.Range("AA9").Value = 0.475
$ws.Range("AB9").Value = -1
$ws.Range("AC9").Value = 1
$ws.Range("B9").Value = 6814330
$ws.Range("G9").Value = "NK Aluminij"
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "H"
$ws.Range("K9").Value = 1.363
$ws.Range("L9").Value = 4.5
$ws.Range("M9").Value = 7
$ws.Range("N9").Value = 1.4
$ws.Range("O9").Value = 4.5
$ws.Range("P9").Value = 7
$ws.Range("Q9").Value = -1.25
$ws.Range("R9").Value = 1.85
$ws.Range("S9").Value = 1.95
$ws.Range("T9").Value = 2.75
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 2
$ws.Range("W9").Value = 0.3999999999999999
$ws.Range("X9").Value = -1
$ws.Range("AA10").Value = 0.4375
$ws.Range("AB10").Value = -0.5
$ws.Range("AC10").Value = 0.425
$ws.Range("B10").Value = 6814328
$ws.Range("G10").Value = "NK Bravo"
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = "D"
$ws.Range("K10").Value = 2.35
$ws.Range("L10").Value = 3.1
$ws.Range("M10").Value = 2.9
$ws.Range("N10").Value = 2.15
$ws.Range("O10").Value = 3.1
$ws.Range("P10").Value = 3.3
$ws.Range("Q10").Value = -0.25
$ws.Range("R10").Value = 1.925
$ws.Range("S10").Value = 1.875
$ws.Range("T10").Value = 2.25
$ws.Range("U10").Value = 1.95
$ws.Range("V10").Value = 1.85
$ws.Range("W10").Value = -1
$ws.Range("X10").Value = 2.1
$ws.Range("G12").Value = "NK Maribor"
$ws.Range("F14").Value = "NK Domzale"
$ws.Range("G20").Value = "NK Domzale"
$ws.Range("F21").Value = "NK Maribor"
$ws.Range("F23").Value = "NK Domzale"
$ws.Range("G23").Value = "NK Maribor"
$ws.Range("F27").Value = "NK Maribor"
$ws.Range("G29").Value = "NK Domzale"
$ws.Range("G33").Value = "NK Domzale"
$ws.Range("G34").Value = "NK Maribor"
$ws.Range("F36").Value = "NK Maribor"
$ws.Range("F40").Value = "NK Domzale"
$ws.Range("F43").Value = "NK Maribor"
$ws.Range("F44").Value = "NK Domzale"
$ws.Range("G46").Value = "NK Maribor"
$ws.Range("G50").Value = "NK Domzale"
$ws.Range("F52").Value = "NK Domzale"
$ws.Range("F53").Value = "NK Maribor"
$ws.Range("F57").Value = "NK Maribor"
$ws.Range("G58").Value = "NK Domzale"
$ws.Range("G61").Value = "NK Maribor"
$ws.Range("F63").Value = "NK Domzale"
$ws.Range("F68").Value = "NK Maribor"
$ws.Range("G68").Value = "NK Domzale"
$ws.Range("G73").Value = "NK Maribor"
$ws.Range("F74").Value = "NK Domzale"
$ws.Range("G76").Value = "NK Maribor"
$ws.Range("G79").Value = "NK Domzale"
$ws.Range("F81").Value = "NK Maribor"
$ws.Range("F83").Value = "NK Domzale"
$ws.Range("G85").Value = "NK Maribor"
$ws.Range("G86").Value = "NK Domzale"
$ws.Range("G92").Value = "NK Maribor"
$ws.Range("G93").Value = "NK Domzale"
$ws.Range("F96").Value = "NK Maribor"
$ws.Range("F99").Value = "NK Domzale"
$ws.Range("G101").Value = "NK Maribor"
$ws.Range("G105").Value = "NK Maribor"
$ws.Range("F107").Value = "NK Domzale"
$ws.Range("G109").Value = "NK Domzale"
$ws.Range("F110").Value = "NK Maribor"
$ws.Range("F115").Value = "NK Domzale"
$ws.Range("G115").Value = "NK Maribor"
$ws.Range("G119").Value = "NK Domzale"
$ws.Range("G124").Value = "NK Maribor"
$ws.Range("F127").Value = "NK Maribor"
$ws.Range("F130").Value = "NK Domzale"
$ws.Range("F133").Value = "NK Domzale"
$ws.Range("F134").Value = "NK Maribor"
$ws.Range("G137").Value = "NK Domzale"
$ws.Range("B138").Value = 6837117
$ws.Range("E138").Value = 45388.63541666666
$ws.Range("F138").Value = "NS Mura"
$ws.Range("G138").Value = "NK Celje"
$ws.Range("K138").Value = 5.25
$ws.Range("L138").Value = 4.2
$ws.Range("M138").Value = 1.5
$ws.Range("N138").Value = 5.25
$ws.Range("O138").Value = 4.2
$ws.Range("P138").Value = 1.5
$ws.Range("R138").Value = 1.975
$ws.Range("S138").Value = 1.825
$ws.Range("U138").Value = 1.95
$ws.Range("V138").Value = 1.85
$ws.Range("B139").Value = 6814434
$ws.Range("E139").Value = 45389.41666666666
$ws.Range("F139").Value = "NK Bravo"
$ws.Range("G139").Value = "NK Maribor"
$ws.Range("K139").Value = 1.833
$ws.Range("M139").Value = 4
$ws.Range("N139").Value = 1.7
$ws.Range("O139").Value = 3.4
$ws.Range("P139").Value = 4.5
$ws.Range("Q139").Value = -0.75
$ws.Range("R139").Value = 1.975
$ws.Range("S139").Value = 1.825
$ws.Range("U139").Value = 1.775
$ws.Range("V139").Value = 2.025

# --- Remove the AA139 cell (no longer populated in the refreshed row) ---
$ws.Range("AA139").ClearContents()

# --- Drop the last two fixture rows (140, 141); the sheet's used range shrinks to A1:AC139 ---
$ws.Rows("140:141").Delete()

Write-Host "edit applied"
